$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename some RefDes labels (simple text swaps) ---
$ws.Range("A2").Value = "3PDT1"
$ws.Range("A5").Value = "LEVEL"
$ws.Range("A6").Value = "DRIVE"
$ws.Range("A8").Value = "LED1"
$ws.Range("A9").Value = "LED1_HLD"
$ws.Range("A10").Value = "FOCUS"

# --- Update row 11 (THK,EDG pot -> BODY,EDGE pot, 5K -> 25K) ---
$ws.Range("A11").Value = "BODY, EDGE"
$ws.Range("B11").Value = "25K OHM Logarithmic Taper Potentiometer Round Knurled Plastic Shaft PCB 9mm SKU: A-1880"
$ws.Range("C11").Value = "RV09AF-40-20K-A25K"

# --- Insert a new table row just above the Totals row (row 13) ---
$ws.Rows("13").Insert()

# Fill in the new data row with the added BOM line item
$ws.Range("A13").Value = "RESONANCE, VOICING"
$ws.Range("B13").Value = "Micro Toggle Switch SPDT On-On SKU: A-3643"
$ws.Range("D13").Formula = '=(LEN(A13)-LEN(SUBSTITUTE(A13,",",""))+1)'
$ws.Range("E13").Value = 0.46
$ws.Range("F13").Formula = "='TAYDA ORDER'!`$E13*'TAYDA ORDER'!`$D13"

# --- Resize the table to include the new row ---
$tbl = $ws.ListObjects("Table_3")
$tbl.Resize($ws.Range("A1:F14"))

# --- Let Excel recompute the autofit column widths for the changed text ---
$ws.Columns("A:B").AutoFit()

# --- Update the active selection (cosmetic) ---
$ws.Range("B6").Select()
